{"js": "// Replace the 25 \"three-digit x one-digit\" multiplication problems in the\n// document's table with their new values, preserving run/paragraph\n// formatting (font, size, alignment) by doing an in-place text replace on\n// the existing range rather than rewriting the whole cell/paragraph body.\nconst replacements = [\n  [\"606\u00d74=2424\", \"429\u00d77=3003\"],\n  [\"582\u00d77=4074\", \"114\u00d73=342\"],\n  [\"912\u00d74=3648\", \"723\u00d79=6507\"],\n  [\"921\u00d73=2763\", \"922\u00d78=7376\"],\n  [\"742\u00d76=4452\", \"278\u00d73=834\"],\n  [\"323\u00d78=2584\", \"336\u00d76=2016\"],\n  [\"559\u00d73=1677\", \"132\u00d79=1188\"],\n  [\"472\u00d78=3776\", \"550\u00d78=4400\"],\n  [\"430\u00d76=2580\", \"936\u00d75=4680\"],\n  [\"554\u00d79=4986\", \"146\u00d78=1168\"],\n  [\"775\u00d74=3100\", \"684\u00d73=2052\"],\n  [\"621\u00d74=2484\", \"278\u00d78=2224\"],\n  [\"522\u00d74=2088\", \"625\u00d72=1250\"],\n  [\"486\u00d73=1458\", \"361\u00d77=2527\"],\n  [\"669\u00d79=6021\", \"743\u00d77=5201\"],\n  [\"935\u00d77=6545\", \"450\u00d73=1350\"],\n  [\"531\u00d74=2124\", \"807\u00d72=1614\"],\n  [\"292\u00d77=2044\", \"888\u00d74=3552\"],\n  [\"527\u00d74=2108\", \"590\u00d73=1770\"],\n  [\"589\u00d75=2945\", \"488\u00d72=976\"],\n  [\"658\u00d76=3948\", \"372\u00d75=1860\"],\n  [\"124\u00d76=744\", \"801\u00d75=4005\"],\n  [\"104\u00d77=728\", \"378\u00d78=3024\"],\n  [\"678\u00d74=2712\", \"581\u00d73=1743\"],\n  [\"847\u00d75=4235\", \"607\u00d74=2428\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" multiplication problems in the\n# document's table with their new values. Using Find/Execute with a\n# Replacement.Text (Find What / Replace With) performs an in-place text\n# substitution that preserves the existing run/paragraph formatting\n# (font, size, alignment) instead of rewriting the whole paragraph.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"606\u00d74=2424\", \"429\u00d77=3003\"),\n    @(\"582\u00d77=4074\", \"114\u00d73=342\"),\n    @(\"912\u00d74=3648\", \"723\u00d79=6507\"),\n    @(\"921\u00d73=2763\", \"922\u00d78=7376\"),\n    @(\"742\u00d76=4452\", \"278\u00d73=834\"),\n    @(\"323\u00d78=2584\", \"336\u00d76=2016\"),\n    @(\"559\u00d73=1677\", \"132\u00d79=1188\"),\n    @(\"472\u00d78=3776\", \"550\u00d78=4400\"),\n    @(\"430\u00d76=2580\", \"936\u00d75=4680\"),\n    @(\"554\u00d79=4986\", \"146\u00d78=1168\"),\n    @(\"775\u00d74=3100\", \"684\u00d73=2052\"),\n    @(\"621\u00d74=2484\", \"278\u00d78=2224\"),\n    @(\"522\u00d74=2088\", \"625\u00d72=1250\"),\n    @(\"486\u00d73=1458\", \"361\u00d77=2527\"),\n    @(\"669\u00d79=6021\", \"743\u00d77=5201\"),\n    @(\"935\u00d77=6545\", \"450\u00d73=1350\"),\n    @(\"531\u00d74=2124\", \"807\u00d72=1614\"),\n    @(\"292\u00d77=2044\", \"888\u00d74=3552\"),\n    @(\"527\u00d74=2108\", \"590\u00d73=1770\"),\n    @(\"589\u00d75=2945\", \"488\u00d72=976\"),\n    @(\"658\u00d76=3948\", \"372\u00d75=1860\"),\n    @(\"124\u00d76=744\", \"801\u00d75=4005\"),\n    @(\"104\u00d77=728\", \"378\u00d78=3024\"),\n    @(\"678\u00d74=2712\", \"581\u00d73=1743\"),\n    @(\"847\u00d75=4235\", \"607\u00d74=2428\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n\n    if (-not $found) {\n        Write-Output \"WARNING: text not found: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
